$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint 4 backlog: reassign reviewers for a handful of user stories.
$ws.Range("F2").Value = "Yasser"

$ws.Range("E3").Value = "Hosam"
$ws.Range("F3").Value = "Mostafa"

$ws.Range("E4").Value = "Farid"
$ws.Range("F4").Value = "Mohamed Ayman"

$ws.Range("E6").Value = "Mark"
$ws.Range("F6").Value = "Hosam"

$ws.Range("F7").Value = "Ahmed El-Sherif"

$ws.Range("E9").Value = "Mostafa"
$ws.Range("F9").Value = "Mark"

$ws.Range("E10").Value = "Ahmed El-Sherif"
